$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = "maa://21246 (91.41), maa://36684 (95.5), ***maa://22731 (6.25)"
$ws.Range("AF2").Value = "maa://25251 (92.31), ***maa://21730 (25.33), ***maa://39501 (16.67), **maa://36675 (50.0)"
$ws.Range("P3").Value = "maa://21249 (94.44), maa://26254 (96.55)"
$ws.Range("T3").Value = "maa://24617 (89.74), **maa://20790 (43.48), ***maa://37170 (16.92), maa://45854 (88.24)"
$ws.Range("AB3").Value = "maa://24390 (94.59)"
$ws.Range("X4").Value = "**maa://32495 (48.7), ***maa://31785 (22.22), maa://43217 (89.29), ***maa://36683 (28.26)"
$ws.Range("AF4").Value = "*maa://30062 (64.0), ***maa://26209 (13.04), *maa://39394 (64.0)"
$ws.Range("AB5").Value = "*maa://29863 (64.86), ***maa://22752 (12.5), **maa://26013 (37.5)"
$ws.Range("P7").Value = "maa://22750 (91.67)"
$ws.Range("A8").Value = "更新日期：2025.02.15 19:34:59"
$ws.Range("X8").Value = "maa://21411 (95.92)"
$ws.Range("X9").Value = "maa://26223 (97.87)"
$ws.Range("AF9").Value = "maa://26206 (88.43), *maa://22865 (50.94)"
$ws.Range("D10").Value = "***maa://25695 (18.72), ***maa://34206 (20.0), ***maa://39951 (15.69), ***maa://39243 (25.0), *maa://45271 (55.17)"
$ws.Range("X10").Value = "maa://22301 (97.75), maa://45828 (90.0), maa://22726 (100.0)"
$ws.Range("AB12").Value = "maa://23669 (95.5), maa://36677 (93.22), maa://39872 (91.67)"
$ws.Range("X13").Value = "maa://34957 (82.43), *maa://22768 (51.61)"
$ws.Range("D16").Value = "maa://21441 (96.4), maa://36679 (94.34), maa://37650 (97.14)"
$ws.Range("T16").Value = "maa://22729 (94.97), *maa://28648 (69.12), maa://36674 (82.69)"
$ws.Range("X16").Value = "maa://28501 (98.02), maa://28051 (96.0)"
$ws.Range("AB16").Value = "maa://26228 (95.92)"
$ws.Range("AF16").Value = "*maa://23911 (65.09), maa://27755 (93.55)"
$ws.Range("D18").Value = "maa://24570 (97.32)"
$ws.Range("L18").Value = "maa://22466 (90.06), *maa://22732 (51.14)"
$ws.Range("D20").Value = "maa://21432 (89.94), maa://25198 (93.58), *maa://20795 (50.77), maa://36680 (91.18)"
$ws.Range("L20").Value = "maa://41331 (85.71)"
$ws.Range("D21").Value = "maa://21261 (97.56)"
$ws.Range("D23").Value = "***maa://28036 (28.77), *maa://41753 (56.25)"
$ws.Range("L23").Value = "maa://39756 (95.55), maa://39875 (94.37)"
$ws.Range("X24").Value = "maa://29988 (84.58), maa://23504 (93.1), **maa://22892 (40.14), *maa://25141 (77.1), *maa://36663 (77.63), ***maa://22815 (23.08)"
$ws.Range("AB25").Value = "maa://31215 (87.61), maa://24516 (80.22), maa://26001 (87.5)"
$ws.Range("D26").Value = "maa://41802 (93.75)"
$ws.Range("H26").Value = "maa://24913 (92.13)"
$ws.Range("AB26").Value = "maa://42235 (94.85)"
$ws.Range("H27").Value = "**maa://21283 (47.37), *maa://39601 (80.0), maa://34494 (97.14), **maa://36665 (50.0)"
$ws.Range("X28").Value = "maa://39929 (90.6), maa://41749 (90.48), ***maa://39723 (13.89)"
$ws.Range("D29").Value = "maa://31694 (98.15)"
$ws.Range("L29").Value = "maa://28432 (93.43), *maa://28440 (79.63), maa://31400 (98.81), *maa://28650 (71.43)"
$ws.Range("T32").Value = "maa://42859 (95.69), maa://41108 (88.0), maa://41238 (97.09), maa://45523 (100.0)"
$ws.Range("L35").Value = "maa://41296 (96.18)"
$ws.Range("L37").Value = "maa://45718 (98.35), maa://45789 (100.0)"
$ws.Range("P39").Value = "maa://24709 (91.45)"
$ws.Range("S39").NumberFormat = "@"
$ws.Range("S39").Value = "3"
$ws.Range("T39").Value = "maa://45788 (82.8), *maa://45790 (75.0), **maa://47079 (50.0)"
$ws.Range("H44").Value = "maa://29768 (98.01), maa://27728 (96.08)"
$ws.Range("H46").Value = "maa://35931 (92.63), maa://43901 (92.0)"
$ws.Range("H47").Value = "maa://27410 (96.44), maa://29661 (97.3), maa://28038 (84.62)"
$ws.Range("H62").Value = "maa://42981 (95.12), maa://43903 (100.0)"
